$wb = $excel.ActiveWorkbook

# --- Sheet 1: nhap-linhkien -> nhap-thanhpham ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "nhap-thanhpham"

# Drop the unused trailing columns (G:J) so the sheet shrinks back to A:F
$ws1.Columns("G:J").Delete()

# New header row
$ws1.Range("A1").Value = "Tên Hàng"
$ws1.Range("B1").Value = "MCU"
$ws1.Range("C1").Value = "Sổ Hợp Đồng"
$ws1.Range("D1").Value = "Chip"
$ws1.Range("E1").Value = "Ngày Nhập"
$ws1.Range("F1").Value = "Số Lượng"

# New data row
$ws1.Range("A2").Value = "RES 0R 5% 3/4W 2010"
$ws1.Range("B2").Value = "mcu01"
$ws1.Range("C2").Value = "006-21/DT-BS"
$ws1.Range("D2").Value = "chip01"
# Force the date-looking string to stay text instead of being parsed as a
# serial date value (leading apostrophe = "enter as text" in Excel), then
# strip the quote-prefix style back off so the cell stays unstyled.
$ws1.Range("E2").Value = "'2021-10-11"
$ws1.Range("E2").Style = "Normal"
$ws1.Range("F2").Value = 10

# --- Sheet 2: xuat-linhkien -> xuat-thanhpham ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "xuat-thanhpham"

$ws2.Columns("G:J").Delete()

$ws2.Range("A1").Value = "Tên Hàng"
$ws2.Range("B1").Value = "MCU"
$ws2.Range("C1").Value = "Sổ Hợp Đồng"
$ws2.Range("D1").Value = "Chip"
$ws2.Range("E1").Value = "Ngày Nhập"
$ws2.Range("F1").Value = "Số Lượng"

# --- Sheet 3: ton-linhkien -> ton-thanhpham ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "ton-thanhpham"

$ws3.Range("A1").Value = "Tên Hàng"
$ws3.Range("B1").Value = "Số Lượng"
$ws3.Range("C1").Value = "Đơn Vị Tính"

$ws3.Range("A2").Value = "RES 0R 5% 3/4W 2010"
$ws3.Range("B2").Value = 10
$ws3.Range("C2").Value = "none"
